$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "1.00", "19.00") are preserved exactly instead of being
# coerced into numbers and losing trailing zeros / precision.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.059.42"
$ws.Range("E2").Value = "  -6.11%  "
$ws.Range("D3").Value = "3.566.22"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "395.68"
$ws.Range("E5").Value = "  -5.52%  "
$ws.Range("D6").Value = "125.79"
$ws.Range("E6").Value = "  -5.55%  "
$ws.Range("D7").Value = "3.554.90"
$ws.Range("E7").Value = "  -1.88%  "
$ws.Range("E8").Value = "  -8.72%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "0.693"
$ws.Range("E10").Value = "  -10.10%  "
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  -15.37%  "
$ws.Range("D12").Value = "0.0000358"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "39.57"
$ws.Range("E13").Value = "  -8.04%  "
$ws.Range("D14").Value = "4.122.00"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "9.33"
$ws.Range("E15").Value = "  -5.84%  "
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").Value = "3.615.72"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "12.96"
$ws.Range("E18").Value = "  +4.86%  "
$ws.Range("D19").Value = "19.00"
$ws.Range("E19").Value = "  -7.09%  "
$ws.Range("D20").Value = "64.459.74"
$ws.Range("E20").Value = "  -5.10%  "
$ws.Range("E21").Value = "  -9.83%  "
$ws.Range("D22").Value = "401.43"
$ws.Range("E22").Value = "  -14.40%  "
$ws.Range("D23").Value = "14.08"
$ws.Range("E23").Value = "  +5.16%  "
$ws.Range("D24").Value = "81.96"
$ws.Range("E24").Value = "  -7.72%  "
$ws.Range("D25").Value = "2.90"
$ws.Range("E25").Value = "  -7.44%  "
$ws.Range("D26").Value = "34.32"
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("D27").Value = "5.34"
$ws.Range("E27").Value = "  +9.61%  "
$ws.Range("E28").Value = "  -9.91%  "
$ws.Range("D29").Value = "8.89"
$ws.Range("E29").Value = "  -11.51%  "
$ws.Range("D30").Value = "12.01"
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").Value = "2.59"
$ws.Range("E31").Value = "  -7.04%  "
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("D33").Value = "6.90"
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").Value = "0.153"
$ws.Range("E34").Value = "  -5.57%  "
$ws.Range("D35").Value = "37.80"
$ws.Range("E35").Value = "  -7.20%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "54.13"
$ws.Range("E37").Value = "  -4.68%  "
$ws.Range("D38").Value = "0.0444"
$ws.Range("E38").Value = "  -10.38%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").Value = "2.78"
$ws.Range("E40").Value = "  +18.53%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0636"
$ws.Range("E41").Value = "  -9.61%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.132"
$ws.Range("E42").Value = "  -9.35%  "
$ws.Range("D43").Value = "141.12"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("E44").Value = "  +13.74%  "
$ws.Range("D45").Value = "4.26"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "25.22"
$ws.Range("E47").Value = "  +16.16%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  -9.69%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "1.97"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "2.50"
$ws.Range("E50").Value = "  -7.38%  "
$ws.Range("E51").Value = "  -9.62%  "

# Restore default styling on column D (remove the temporary text
# number format) now that values are safely written as text.
$ws.Range("D2:D51").Style = "Normal"
